# Duplicate the existing 26 data rows (rows 2-27) into rows 28-53,
# continuing the "id" sequence from 27 to 52.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Metal Mario - Tennis", "Metal Mario", "Mario Sports Superstars", "Mario Sports Superstars", "Card"),
    @("Super Mario Cereal", "Mario Cereal", "Others", "Kellogs", "Card"),
    @("Baby Mario - Soccer", "Baby Mario", "Mario Sports Superstars", "Mario Sports Superstars", "Card"),
    @("Metal Mario - Soccer", "Metal Mario", "Mario Sports Superstars", "Mario Sports Superstars", "Card"),
    @("Mario - Soccer", "Mario", "Mario Sports Superstars", "Mario Sports Superstars", "Card"),
    @("Mario", "Mario", "Super Mario Bros.", "Super Mario", "Figure"),
    @("8-Bit Mario Modern Color", "Mario", "8-bit Mario", "Super Mario", "Figure"),
    @("Dr. Mario", "Mario", "Super Smash Bros.", "Super Mario", "Figure"),
    @("Baby Mario - Horse Racing", "Baby Mario", "Mario Sports Superstars", "Mario Sports Superstars", "Card"),
    @("Mario - Cat", "Mario", "Super Mario Bros.", "Super Mario", "Figure"),
    @("Baby Mario - Golf", "Baby Mario", "Mario Sports Superstars", "Mario Sports Superstars", "Card"),
    @("Mario - Wedding", "Mario", "Super Mario Bros.", "Super Mario", "Figure"),
    @("Metal Mario - Golf", "Metal Mario", "Mario Sports Superstars", "Mario Sports Superstars", "Card"),
    @("Mario", "Mario", "Super Smash Bros.", "Super Mario", "Figure"),
    @("Mario - Tennis", "Mario", "Mario Sports Superstars", "Mario Sports Superstars", "Card"),
    @("Baby Mario - Tennis", "Baby Mario", "Mario Sports Superstars", "Mario Sports Superstars", "Card"),
    @("Mario - Gold Edition", "Mario", "Super Mario Bros.", "Super Mario", "Figure"),
    @("Mario - Power Up Band", "Mario", "Super Nintendo World", "Super Mario", "Band"),
    @("Metal Mario - Horse Racing", "Metal Mario", "Mario Sports Superstars", "Mario Sports Superstars", "Card"),
    @("Baby Mario - Baseball", "Baby Mario", "Mario Sports Superstars", "Mario Sports Superstars", "Card"),
    @("Mario - Golf", "Mario", "Mario Sports Superstars", "Mario Sports Superstars", "Card"),
    @("Mario - Silver Edition", "Mario", "Super Mario Bros.", "Super Mario", "Figure"),
    @("Mario - Horse Racing", "Mario", "Mario Sports Superstars", "Mario Sports Superstars", "Card"),
    @("8-Bit Mario Classic Color", "Mario", "8-bit Mario", "Super Mario", "Figure"),
    @("Metal Mario - Baseball", "Metal Mario", "Mario Sports Superstars", "Mario Sports Superstars", "Card"),
    @("Mario - Baseball", "Mario", "Mario Sports Superstars", "Mario Sports Superstars", "Card")
)

$startRow = 28
$id = 27

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $rec = $data[$i]

    $ws.Cells.Item($row, 1).Value = $id
    $ws.Cells.Item($row, 2).Value = $rec[0]
    $ws.Cells.Item($row, 3).Value = $rec[1]
    $ws.Cells.Item($row, 4).Value = $rec[2]
    $ws.Cells.Item($row, 5).Value = $rec[3]
    $ws.Cells.Item($row, 6).Value = $rec[4]

    $id = $id + 1
}
